$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.684.20"
$ws.Range("E2").Value = "  -2.42%  "

$ws.Range("D3").Value = "3.395.39"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'406.70"
$ws.Range("E5").Value = "  -2.26%  "

$ws.Range("D6").Value = "'133.50"
$ws.Range("E6").Value = "  +7.92%  "

$ws.Range("D7").Value = "'0.591"
$ws.Range("E7").Value = "  -2.32%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.670"
$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").Value = "'0.121"
$ws.Range("E10").Value = "  -7.14%  "

$ws.Range("D11").Value = "'42.59"
$ws.Range("E11").Value = "  +2.84%  "

$ws.Range("E12").Value = "  -1.13%  "

$ws.Range("D13").Value = "3.922.28"
$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").Value = "'8.41"
$ws.Range("E14").Value = "  -2.03%  "

$ws.Range("D15").Value = "'19.76"
$ws.Range("E15").Value = "  -0.95%  "

$ws.Range("D16").Value = "3.384.26"
$ws.Range("E16").Value = "  -2.46%  "

$ws.Range("D17").Value = "61.635.64"
$ws.Range("E17").Value = "  -2.27%  "

$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").Value = "'10.99"
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("E20").Value = "  -8.66%  "

$ws.Range("D21").Value = "'3.20"
$ws.Range("E21").Value = "  -3.99%  "

$ws.Range("D22").Value = "'85.15"
$ws.Range("E22").Value = "  +2.73%  "

$ws.Range("D23").Value = "'315.90"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").Value = "'12.82"

$ws.Range("D25").Value = "'3.14"
$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("E26").Value = "  +11.23%  "

$ws.Range("D27").Value = "'8.34"
$ws.Range("E27").Value = "  +5.84%  "

$ws.Range("D28").Value = "'29.56"
$ws.Range("E28").Value = "  -4.65%  "

$ws.Range("D29").Value = "'7.59"
$ws.Range("E29").Value = "  -2.16%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.69"
$ws.Range("E30").Value = "  +4.90%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.116"
$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("E32").Value = "  -2.56%  "

$ws.Range("D33").Value = "'11.39"
$ws.Range("E33").Value = "  -1.97%  "

$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").Value = "'41.17"
$ws.Range("E35").Value = "  -2.18%  "

$ws.Range("D36").Value = "'0.0482"
$ws.Range("E36").Value = "  -2.08%  "

$ws.Range("D37").Value = "'51.71"
$ws.Range("E37").Value = "  -1.05%  "

$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  -3.21%  "

$ws.Range("D41").Value = "'140.38"
$ws.Range("E41").Value = "  +3.84%  "

$ws.Range("D42").Value = "'1.98"
$ws.Range("E42").Value = "  -1.29%  "

$ws.Range("E43").Value = "  +4.73%  "

$ws.Range("E44").Value = "  -2.02%  "

$ws.Range("D45").Value = "'4.00"
$ws.Range("E45").Value = "  +2.35%  "

$ws.Range("D46").Value = "'16.60"
$ws.Range("E46").Value = "  -3.44%  "

$ws.Range("E47").Value = "  -0.86%  "

$ws.Range("D48").Value = "'21.35"
$ws.Range("E48").Value = "  -3.11%  "

$ws.Range("D49").Value = "2.114.99"
$ws.Range("E49").Value = "  -3.17%  "

$ws.Range("E50").Value = "  -5.02%  "

$ws.Range("D51").Value = "'1.90"
$ws.Range("E51").Value = "  -0.05%  "
